$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 59.4
$ws.Range("I11").Value = 59.4
$ws.Range("K11").Value = 59.4
$ws.Range("M11").Value = 80.59999999999999
$ws.Range("I19").Value = 2000
$ws.Range("K19").Value = 2000
$ws.Range("M19").Value = -1825
$ws.Range("H40").Value = 5452.1304
$ws.Range("I40").Value = 985.5714
$ws.Range("J40").Value = 7406.25
$ws.Range("K40").Value = 985.5714
$ws.Range("L40").Value = 7406.25
$ws.Range("M40").Value = -810.5714
$ws.Range("N40").Value = -7756.25
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H96").Value = 2400.818
$ws.Range("I96").Value = 2815.5715
$ws.Range("J96").Value = 1675
$ws.Range("K96").Value = 8446.7145
$ws.Range("L96").Value = 5025
$ws.Range("M96").Value = -7073.7145
$ws.Range("N96").Value = -7771
$ws.Range("H106").Value = 3999.5
$ws.Range("I106").Value = 3999.5
$ws.Range("K106").Value = 3999.5
$ws.Range("M106").Value = -3368.5
$ws.Range("H107").Value = 1498.091
$ws.Range("I107").Value = 886.6667
$ws.Range("K107").Value = 886.6667
$ws.Range("M107").Value = 1033.3333

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 300
$ws.Range("K5").Value = 300
$ws.Range("M5").Value = -188
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H122").Value = 7386.3887
$ws.Range("I122").Value = 7597.933
$ws.Range("K122").Value = 22793.799
$ws.Range("M122").Value = -20343.799
$ws.Range("H132").Value = 3613.6
$ws.Range("J132").Value = 5422.5
$ws.Range("L132").Value = 16267.5
$ws.Range("N132").Value = -21327.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -185
$ws.Range("H134").Value = 6666.1113
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = -465

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5840.231
$ws.Range("I31").Value = 5885.4
$ws.Range("J31").Value = 5812
$ws.Range("K31").Value = 5885.4
$ws.Range("L31").Value = 5812
$ws.Range("M31").Value = -5590.4
$ws.Range("N31").Value = -6402
$ws.Range("H34").Value = 5840.231
$ws.Range("I34").Value = 5885.4
$ws.Range("J34").Value = 5812
$ws.Range("K34").Value = 5885.4
$ws.Range("L34").Value = 5812
$ws.Range("M34").Value = -5683.4
$ws.Range("N34").Value = -6216
$ws.Range("H86").Value = 4729.2
$ws.Range("I86").Value = 4729.2
$ws.Range("K86").Value = 4729.2
$ws.Range("M86").Value = -3606.2
$ws.Range("H89").Value = 4729.2
$ws.Range("I89").Value = 4729.2
$ws.Range("K89").Value = 23646
$ws.Range("M89").Value = -18030
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 4699.2
$ws.Range("I122").Value = 6298.5
$ws.Range("J122").Value = 3633
$ws.Range("K122").Value = 18895.5
$ws.Range("L122").Value = 10899
$ws.Range("M122").Value = -16445.5
$ws.Range("N122").Value = -15799
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("N127").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 80.5
$ws.Range("I2").Value = 42.5
$ws.Range("J2").Value = 118.5
$ws.Range("K2").Value = 255
$ws.Range("L2").Value = 711
$ws.Range("M2").Value = -142
$ws.Range("N2").Value = -937
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H80").Value = 5145.273
$ws.Range("I80").Value = 4174.75
$ws.Range("J80").Value = 5699.857
$ws.Range("K80").Value = 12524.25
$ws.Range("L80").Value = 17099.571
$ws.Range("M80").Value = -11588.25
$ws.Range("N80").Value = -18971.571
$ws.Range("H83").Value = 5145.273
$ws.Range("I83").Value = 4174.75
$ws.Range("J83").Value = 5699.857
$ws.Range("K83").Value = 37572.75
$ws.Range("L83").Value = 51298.713
$ws.Range("M83").Value = -32892.75
$ws.Range("N83").Value = -60658.713
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -8372
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -29856
$ws.Range("H114").Value = 1593
$ws.Range("I114").Value = 640
$ws.Range("J114").Value = 3499
$ws.Range("K114").Value = 1920
$ws.Range("L114").Value = 10497
$ws.Range("M114").Value = 1334
$ws.Range("N114").Value = -17005
$ws.Range("H117").Value = 1833.3334
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 1833.3334
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -12384.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1189.6
$ws.Range("I102").Value = 1189.6
$ws.Range("K102").Value = 1189.6
$ws.Range("M102").Value = 432.4000000000001
$ws.Range("H122").Value = 1581
$ws.Range("I122").Value = 1302
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 3906
$ws.Range("L122").Value = 5998.5
$ws.Range("M122").Value = -1456
$ws.Range("N122").Value = -10898.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5500
$ws.Range("I7").Value = 5500
$ws.Range("K7").Value = 5500
$ws.Range("M7").Value = -5388
$ws.Range("H40").Value = 6215.091
$ws.Range("I40").Value = 6040.6665
$ws.Range("K40").Value = 6040.6665
$ws.Range("M40").Value = -5904.6665
$ws.Range("H46").Value = 4268
$ws.Range("J46").Value = 4268
$ws.Range("L46").Value = 4268
$ws.Range("N46").Value = -4644
$ws.Range("H126").Value = 5500
$ws.Range("I126").Value = 5500
$ws.Range("K126").Value = 16500
$ws.Range("M126").Value = -14030
$ws.Range("H132").Value = 6225
$ws.Range("I132").Value = 6370
$ws.Range("K132").Value = 19110
$ws.Range("M132").Value = -16580

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4499.75
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4499.75
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5747.75
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4499.75
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 22498.75
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -28738.75
$ws.Range("H81").Value = 50000.75
$ws.Range("I81").Value = 50000.75
$ws.Range("K81").Value = 100001.5
$ws.Range("M81").Value = -98940.5
$ws.Range("H84").Value = 50000.75
$ws.Range("I84").Value = 50000.75
$ws.Range("K84").Value = 500007.5
$ws.Range("M84").Value = -494703.5
$ws.Range("H126").Value = 2071.5
$ws.Range("I126").Value = 2071.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6214.5
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 6519.222
$ws.Range("I132").Value = 5174.7144
$ws.Range("K132").Value = 15524.1432
$ws.Range("M132").Value = -12994.1432
